$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-2d/'' target=''_blank''>/Aura Designer 2D</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura/'' target=''_blank''>/Aura</a> 3. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-prime/'' target=''_blank''>/Aura Prime</a> 4. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/avamour/'' target=''_blank''>/Avamour</a> 5. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/avancer-e-sense-2/'' target=''_blank''>/Avancer E-Sense</a> 6. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/caliber-hi-speed-metallic/'' target=''_blank''>/Caliber Hi Speed Metallic</a> 7. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/decair/'' target=''_blank''>/Dec''Air</a> 8. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/dorado/'' target=''_blank''>/Dorado</a> 9. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/flyleaf-2/'' target=''_blank''>/Flyleaf</a> 10. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/high-speed-decora/'' target=''_blank''>/High Speed Decora</a> 11. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/hs-decora-scion/'' target=''_blank''>/HS Decora Scion</a> 12. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/jura/'' target=''_blank''>/Jura</a>'
$ws.Range("B16").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/triton-electroplated/'' target=''_blank''>/Trinton Electroplate</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/trigger-neo/'' target=''_blank''>/Trigger Neo</a>'
$ws.Range("B17").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-2d/'' target=''_blank''>/Aura Designer 2D</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura/'' target=''_blank''>/Aura</a> 3. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-prime/'' target=''_blank''>/Aura Prime</a> 4. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/avamour/'' target=''_blank''>/Avamour</a> 5. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/avancer-e-sense-2/'' target=''_blank''>/Avancer E-Sense</a> 6. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/caliber-hi-speed-metallic/'' target=''_blank''>/Caliber Hi Speed Metallic</a> 7. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/decair/'' target=''_blank''>/Dec''Air</a> 8. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/dorado/'' target=''_blank''>/Dorado</a> 9. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/flyleaf-2/'' target=''_blank''>/Flyleaf</a> 10. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/high-speed-decora/'' target=''_blank''>/High Speed Decora</a> 11. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/hs-decora-scion/'' target=''_blank''>/HS Decora Scion</a> 12. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/jura/'' target=''_blank''>/Jura</a>  13. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/senoprime-2/'' target=''_blank''>/Senoprime</a> 14. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/splitz/'' target=''_blank''>/Splitz</a> 15. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-3d/'' target=''_blank''>/Aura Designer 3D</a>'
$ws.Range("B18").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-2d/'' target=''_blank''>/Aura Designer 2D</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-prime/'' target=''_blank''>/Aura Prime</a> 3. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-3d/'' target=''_blank''>/Aura Designer 3D</a>'
$ws.Range("B19").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-designer-2d/'' target=''_blank''>/Aura Designer 2D</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-prime/'' target=''_blank''>/Aura Prime</a> 3. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/dignita-antidust/'' target=''_blank''>/Dignita Antidust</a> 4. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/jura-antidust/'' target=''_blank''>/Jura Antidust</a> 5. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/avancer-prime-antidust/'' target=''_blank''>/Avancer Prime Antidust</a> 6. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-prime-antidust/'' target=''_blank''>/Aura Prime Antidust</a>'
$ws.Range("B20").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/air-360/'' target=''_blank''>/Air 360</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/air-360-deco/'' target=''_blank''>/Air 360 Deco</a>'
$ws.Range("B21").Value = 'Great! You are one step away. There are many models. Please select the model from the following: 1. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/aura-plus-2/'' target=''_blank''>Aura Plus</a> 2. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/asset-50/'' target=''_blank''>Asset 50</a> 3. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/empower-50/'' target=''_blank''>Empower 50</a> 4. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/entrust-50/'' target=''_blank''>Entrust 50</a> 5. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/high-speed-plus/'' target=''_blank''>High Speed Plus</a> 6. <a href = ''https://www.crompton.co.in/product/consumer-fans/ceiling-fans/hs-decor-plus-2/'' target=''_blank''>High Decor Plus</a>'

$ws.Range("B15").Select()
